$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title text updates (October -> November) ---
$ws.Range("A2").Value = "by State, by Sector, November 2016 and 2015 (Thousand Megawatthours)"
foreach ($addr in @("B6","E6","G6","I6","K6")) { $ws.Range($addr).Value = "November 2016" }
foreach ($addr in @("C6","F6","H6","J6","L6")) { $ws.Range($addr).Value = "November 2015" }

# --- Data cell updates ---
$ws.Range("B7").Value = 840
$ws.Range("C7").Value = 906
$ws.Range("D7").Value = -0.073
$ws.Range("E7").Value = 74
$ws.Range("F7").Value = 75
$ws.Range("G7").Value = 653
$ws.Range("H7").Value = 728
$ws.Range("I7").Value = 13
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 99
$ws.Range("L7").Value = 88

$ws.Range("B8").Value = 73
$ws.Range("C8").Value = 65
$ws.Range("D8").Value = 0.123
$ws.Range("G8").Value = 73
$ws.Range("H8").Value = 65

$ws.Range("B9").Value = 350
$ws.Range("C9").Value = 423
$ws.Range("D9").Value = -0.172
$ws.Range("G9").Value = 244
$ws.Range("H9").Value = 327
$ws.Range("K9").Value = 99
$ws.Range("L9").Value = 88

$ws.Range("B10").Value = 158
$ws.Range("C10").Value = 155
$ws.Range("D10").Value = 0.016
$ws.Range("G10").Value = 148
$ws.Range("H10").Value = 144
$ws.Range("I10").Value = 3
$ws.Range("L10").Value = 0.3

$ws.Range("B11").Value = 175
$ws.Range("C11").Value = 169
$ws.Range("D11").Value = 0.032
$ws.Range("E11").Value = 29
$ws.Range("F11").Value = 23
$ws.Range("G11").Value = 143
$ws.Range("H11").Value = 143

$ws.Range("B12").Value = 22
$ws.Range("C12").Value = 18
$ws.Range("D12").Value = 0.176
$ws.Range("G12").Value = 21
$ws.Range("H12").Value = 18

$ws.Range("C13").Value = 75
$ws.Range("D13").Value = -0.168
$ws.Range("E13").Value = 38
$ws.Range("F13").Value = 44
$ws.Range("G13").Value = 25
$ws.Range("J13").Value = 0.17

$ws.Range("B14").Value = 1279
$ws.Range("C14").Value = 1318
$ws.Range("D14").Value = -0.03
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 1167
$ws.Range("H14").Value = 1188
$ws.Range("I14").Value = 45
$ws.Range("J14").Value = 56
$ws.Range("K14").Value = 60
$ws.Range("L14").Value = 70

$ws.Range("B15").Value = 140
$ws.Range("C15").Value = 123
$ws.Range("D15").Value = 0.139
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 109
$ws.Range("H15").Value = 97
$ws.Range("I15").Value = 24
$ws.Range("J15").Value = 22
$ws.Range("L15").Value = 0.08

$ws.Range("B16").Value = 575
$ws.Range("C16").Value = 613
$ws.Range("D16").Value = -0.062
$ws.Range("G16").Value = 542
$ws.Range("H16").Value = 577
$ws.Range("I16").Value = 17
$ws.Range("L16").Value = 17

$ws.Range("B17").Value = 564
$ws.Range("C17").Value = 582
$ws.Range("D17").Value = -0.031
$ws.Range("G17").Value = 516
$ws.Range("H17").Value = 514
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 15
$ws.Range("K17").Value = 44
$ws.Range("L17").Value = 53

$ws.Range("B18").Value = 2801
$ws.Range("C18").Value = 3291
$ws.Range("D18").Value = -0.149
$ws.Range("E18").Value = 366
$ws.Range("F18").Value = 410
$ws.Range("G18").Value = 2288
$ws.Range("H18").Value = 2718
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 24
$ws.Range("K18").Value = 128
$ws.Range("L18").Value = 139

$ws.Range("B19").Value = 1121
$ws.Range("C19").Value = 1426
$ws.Range("D19").Value = -0.214
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 1117
$ws.Range("H19").Value = 1421

$ws.Range("B20").Value = 533
$ws.Range("C20").Value = 603
$ws.Range("D20").Value = -0.116
$ws.Range("E20").Value = 30
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 497
$ws.Range("H20").Value = 569

$ws.Range("B21").Value = 693
$ws.Range("C21").Value = 735
$ws.Range("D21").Value = -0.057
$ws.Range("E21").Value = 204
$ws.Range("G21").Value = 419
$ws.Range("H21").Value = 442
$ws.Range("K21").Value = 57

$ws.Range("B22").Value = 189
$ws.Range("C22").Value = 211
$ws.Range("D22").Value = -0.102
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 161
$ws.Range("H22").Value = 177
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 30

$ws.Range("B23").Value = 264
$ws.Range("C23").Value = 316
$ws.Range("D23").Value = -0.166
$ws.Range("E23").Value = 127
$ws.Range("F23").Value = 158
$ws.Range("G23").Value = 94
$ws.Range("H23").Value = 108
$ws.Range("K23").Value = 40
$ws.Range("L23").Value = 48

$ws.Range("B24").Value = 5985
$ws.Range("C24").Value = 5800
$ws.Range("D24").Value = 0.032
$ws.Range("E24").Value = 2153
$ws.Range("F24").Value = 2035
$ws.Range("G24").Value = 3771
$ws.Range("H24").Value = 3692
$ws.Range("I24").Value = 9
$ws.Range("J24").Value = 9
$ws.Range("K24").Value = 52
$ws.Range("L24").Value = 63

$ws.Range("B25").Value = 2065
$ws.Range("C25").Value = 2066
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 1344
$ws.Range("F25").Value = 1278
$ws.Range("G25").Value = 714
$ws.Range("H25").Value = 780
$ws.Range("J25").Value = 2

$ws.Range("B26").Value = 1227
$ws.Range("C26").Value = 1158
$ws.Range("D26").Value = 0.06
$ws.Range("E26").Value = 64
$ws.Range("F26").Value = 87
$ws.Range("G26").Value = 1163
$ws.Range("H26").Value = 1071
$ws.Range("K26").Value = 0.02
$ws.Range("L26").Value = 0

$ws.Range("B27").Value = 1164
$ws.Range("C27").Value = 1279
$ws.Range("D27").Value = -0.09
$ws.Range("E27").Value = 268
$ws.Range("F27").Value = 308
$ws.Range("G27").Value = 845
$ws.Range("H27").Value = 909
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 46
$ws.Range("L27").Value = 57

$ws.Range("B28").Value = 118
$ws.Range("C28").Value = 137
$ws.Range("D28").Value = -0.138
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 114
$ws.Range("H28").Value = 134
$ws.Range("I28").Value = "NM"
$ws.Range("J28").Value = 0.11
$ws.Range("L28").Value = 0.15

$ws.Range("B29").Value = 340
$ws.Range("C29").Value = 325
$ws.Range("D29").Value = 0.047
$ws.Range("E29").Value = 26
$ws.Range("F29").Value = 27
$ws.Range("G29").Value = 314
$ws.Range("H29").Value = 296
$ws.Range("J29").Value = 2

$ws.Range("B30").Value = 779
$ws.Range("C30").Value = 584
$ws.Range("D30").Value = 0.334
$ws.Range("E30").Value = 355
$ws.Range("F30").Value = 243
$ws.Range("G30").Value = 424
$ws.Range("H30").Value = 340
$ws.Range("K30").Value = 0.13
$ws.Range("L30").Value = 0.34

$ws.Range("B31").Value = 292
$ws.Range("C31").Value = 251
$ws.Range("D31").Value = 0.16
$ws.Range("E31").Value = 94
$ws.Range("G31").Value = 198
$ws.Range("H31").Value = 162

$ws.Range("B32").Value = 2067
$ws.Range("C32").Value = 1926
$ws.Range("D32").Value = 0.073
$ws.Range("E32").Value = 122
$ws.Range("F32").Value = 145
$ws.Range("G32").Value = 1028
$ws.Range("H32").Value = 873
$ws.Range("I32").Value = 38
$ws.Range("J32").Value = 42
$ws.Range("K32").Value = 880
$ws.Range("L32").Value = 865

$ws.Range("B33").Value = 10
$ws.Range("C33").Value = 10
$ws.Range("D33").Value = 0.038
$ws.Range("F33").Value = 0.32
$ws.Range("G33").Value = 8

$ws.Range("B35").Value = 403
$ws.Range("C35").Value = 411
$ws.Range("D35").Value = -0.019
$ws.Range("E35").Value = 20
$ws.Range("F35").Value = 16
$ws.Range("G35").Value = 205
$ws.Range("H35").Value = 221
$ws.Range("K35").Value = 175
$ws.Range("L35").Value = 170

$ws.Range("B36").Value = 499
$ws.Range("C36").Value = 408
$ws.Range("D36").Value = 0.224
$ws.Range("E36").Value = 14
$ws.Range("F36").Value = 0.05
$ws.Range("G36").Value = 144
$ws.Range("H36").Value = 72
$ws.Range("K36").Value = 341
$ws.Range("L36").Value = 335

$ws.Range("B37").Value = 115
$ws.Range("C37").Value = 103
$ws.Range("D37").Value = 0.113
$ws.Range("G37").Value = 101
$ws.Range("H37").Value = 94
$ws.Range("L37").Value = 7

$ws.Range("B38").Value = 433
$ws.Range("C38").Value = 341
$ws.Range("D38").Value = 0.27
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 300
$ws.Range("H38").Value = 210
$ws.Range("I38").Value = 13
$ws.Range("J38").Value = 16
$ws.Range("K38").Value = 110
$ws.Range("L38").Value = 115

$ws.Range("B39").Value = 184
$ws.Range("C39").Value = 170
$ws.Range("D39").Value = 0.082
$ws.Range("E39").Value = 38
$ws.Range("F39").Value = 42
$ws.Range("G39").Value = 30
$ws.Range("H39").Value = 31
$ws.Range("K39").Value = 116
$ws.Range("L39").Value = 97

$ws.Range("B40").Value = 286
$ws.Range("C40").Value = 331
$ws.Range("D40").Value = -0.135
$ws.Range("E40").Value = 38
$ws.Range("F40").Value = 86
$ws.Range("G40").Value = 104
$ws.Range("H40").Value = 86
$ws.Range("I40").Value = 18
$ws.Range("J40").Value = 20
$ws.Range("K40").Value = 126
$ws.Range("L40").Value = 139

$ws.Range("B41").Value = 136
$ws.Range("C41").Value = 147
$ws.Range("D41").Value = -0.076
$ws.Range("G41").Value = 136
$ws.Range("H41").Value = 147

$ws.Range("B42").Value = 547
$ws.Range("C42").Value = 546
$ws.Range("D42").Value = 0.002
$ws.Range("E42").Value = 8
$ws.Range("F42").Value = 9
$ws.Range("G42").Value = 52
$ws.Range("H42").Value = 46
$ws.Range("J42").Value = 0.19
$ws.Range("K42").Value = 487
$ws.Range("L42").Value = 491

$ws.Range("B43").Value = 290
$ws.Range("C43").Value = 273
$ws.Range("D43").Value = 0.065
$ws.Range("G43").Value = 35
$ws.Range("K43").Value = 255
$ws.Range("L43").Value = 248

$ws.Range("B44").Value = 41
$ws.Range("C44").Value = 46
$ws.Range("D44").Value = -0.1
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 9
$ws.Range("K44").Value = 33
$ws.Range("L44").Value = 36

$ws.Range("B45").Value = 122
$ws.Range("C45").Value = 133
$ws.Range("D45").Value = -0.085
$ws.Range("H45").Value = 3
$ws.Range("K45").Value = 121
$ws.Range("L45").Value = 131

$ws.Range("B46").Value = 93
$ws.Range("C46").Value = 94
$ws.Range("D46").Value = -0.008
$ws.Range("G46").Value = 15
$ws.Range("H46").Value = 18
$ws.Range("J46").Value = 0.19
$ws.Range("L46").Value = 76

$ws.Range("B47").Value = 6809
$ws.Range("C47").Value = 6821
$ws.Range("D47").Value = -0.002
$ws.Range("E47").Value = 153
$ws.Range("F47").Value = 164
$ws.Range("G47").Value = 6219
$ws.Range("H47").Value = 6235
$ws.Range("I47").Value = 6
$ws.Range("J47").Value = 7
$ws.Range("K47").Value = 431
$ws.Range("L47").Value = 415

$ws.Range("B48").Value = 124
$ws.Range("C48").Value = 127
$ws.Range("D48").Value = -0.03
$ws.Range("G48").Value = 9
$ws.Range("H48").Value = 6
$ws.Range("J48").Value = 0.34
$ws.Range("K48").Value = 115
$ws.Range("L48").Value = 121

$ws.Range("B49").Value = 230
$ws.Range("C49").Value = 225
$ws.Range("D49").Value = 0.024
$ws.Range("K49").Value = 223
$ws.Range("L49").Value = 217

$ws.Range("B50").Value = 1725
$ws.Range("C50").Value = 1617
$ws.Range("D50").Value = 0.066
$ws.Range("E50").Value = 138
$ws.Range("F50").Value = 142
$ws.Range("G50").Value = 1561
$ws.Range("H50").Value = 1446
$ws.Range("K50").Value = 25
$ws.Range("L50").Value = 29

$ws.Range("B51").Value = 4730
$ws.Range("C51").Value = 4852
$ws.Range("D51").Value = -0.025
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = 22
$ws.Range("G51").Value = 4641
$ws.Range("H51").Value = 4776
$ws.Range("I51").Value = "NM"
$ws.Range("J51").Value = 6
$ws.Range("K51").Value = 68
$ws.Range("L51").Value = 48

$ws.Range("B52").Value = 3288
$ws.Range("C52").Value = 2782
$ws.Range("D52").Value = 0.182
$ws.Range("E52").Value = 289
$ws.Range("F52").Value = 320
$ws.Range("G52").Value = 2957
$ws.Range("H52").Value = 2419
$ws.Range("I52").Value = 9
$ws.Range("K52").Value = 33
$ws.Range("L52").Value = 35

$ws.Range("B53").Value = 306
$ws.Range("C53").Value = 322
$ws.Range("D53").Value = -0.048
$ws.Range("E53").Value = 33
$ws.Range("F53").Value = 38
$ws.Range("G53").Value = 272
$ws.Range("H53").Value = 282
$ws.Range("I53").Value = "NM"

$ws.Range("B54").Value = 985
$ws.Range("C54").Value = 682
$ws.Range("D54").Value = 0.444
$ws.Range("E54").Value = 24
$ws.Range("F54").Value = 14
$ws.Range("G54").Value = 960
$ws.Range("H54").Value = 667
$ws.Range("I54").Value = "NM"
$ws.Range("L54").Value = 0.24

$ws.Range("B55").Value = 278
$ws.Range("C55").Value = 284
$ws.Range("D55").Value = -0.021
$ws.Range("G55").Value = 245
$ws.Range("H55").Value = 249
$ws.Range("K55").Value = 31
$ws.Range("L55").Value = 33

$ws.Range("B56").Value = 205
$ws.Range("C56").Value = 211
$ws.Range("D56").Value = -0.027
$ws.Range("E56").Value = 21
$ws.Range("F56").Value = 21
$ws.Range("G56").Value = 183
$ws.Range("H56").Value = 187

$ws.Range("B57").Value = 532
$ws.Range("C57").Value = 474
$ws.Range("D57").Value = 0.123
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = 525
$ws.Range("H57").Value = 469
$ws.Range("I57").Value = 5
$ws.Range("L57").Value = 0.18

$ws.Range("B58").Value = 374
$ws.Range("C58").Value = 267
$ws.Range("D58").Value = 0.397
$ws.Range("E58").Value = 16
$ws.Range("F58").Value = 10
$ws.Range("G58").Value = 358
$ws.Range("H58").Value = 257
$ws.Range("J58").Value = 0.3

$ws.Range("B59").Value = 227
$ws.Range("C59").Value = 98
$ws.Range("D59").Value = 1.322
$ws.Range("E59").Value = 23
$ws.Range("F59").Value = 21
$ws.Range("G59").Value = 204
$ws.Range("H59").Value = 76

$ws.Range("B60").Value = 380
$ws.Range("C60").Value = 445
$ws.Range("D60").Value = -0.145
$ws.Range("E60").Value = 170
$ws.Range("F60").Value = 213
$ws.Range("G60").Value = 210
$ws.Range("H60").Value = 232

$ws.Range("B61").Value = 4761
$ws.Range("C61").Value = 4533
$ws.Range("D61").Value = 0.05
$ws.Range("E61").Value = 585
$ws.Range("F61").Value = 595
$ws.Range("G61").Value = 3892
$ws.Range("H61").Value = 3649
$ws.Range("I61").Value = 76
$ws.Range("J61").Value = 81
$ws.Range("K61").Value = 207
$ws.Range("L61").Value = 208

$ws.Range("B62").Value = 3428
$ws.Range("C62").Value = 3108
$ws.Range("D62").Value = 0.103
$ws.Range("E62").Value = 151
$ws.Range("F62").Value = 138
$ws.Range("G62").Value = 3159
$ws.Range("H62").Value = 2857
$ws.Range("I62").Value = 73
$ws.Range("J62").Value = 77
$ws.Range("K62").Value = 45
$ws.Range("L62").Value = 36

$ws.Range("B63").Value = 557
$ws.Range("C63").Value = 645
$ws.Range("D63").Value = -0.137
$ws.Range("E63").Value = 51
$ws.Range("F63").Value = 89
$ws.Range("G63").Value = 462
$ws.Range("H63").Value = 503
$ws.Range("K63").Value = 43
$ws.Range("L63").Value = 51

$ws.Range("B64").Value = 776
$ws.Range("C64").Value = 779
$ws.Range("D64").Value = -0.005
$ws.Range("E64").Value = 383
$ws.Range("F64").Value = 368
$ws.Range("G64").Value = 272
$ws.Range("H64").Value = 289
$ws.Range("K64").Value = 119
$ws.Range("L64").Value = 121

$ws.Range("B65").Value = 141
$ws.Range("C65").Value = 142
$ws.Range("D65").Value = -0.003
$ws.Range("E65").Value = 18
$ws.Range("F65").Value = 17
$ws.Range("G65").Value = 98
$ws.Range("H65").Value = 97
$ws.Range("I65").Value = 20
$ws.Range("J65").Value = 19
$ws.Range("K65").Value = 6

$ws.Range("B66").Value = 19
$ws.Range("C66").Value = 23
$ws.Range("D66").Value = -0.181
$ws.Range("E66").Value = 10
$ws.Range("F66").Value = 12
$ws.Range("H66").Value = 6
$ws.Range("L66").Value = 0.39

$ws.Range("B67").Value = 122
$ws.Range("C67").Value = 119
$ws.Range("D67").Value = 0.031
$ws.Range("E67").Value = 8
$ws.Range("F67").Value = 5
$ws.Range("G67").Value = 93
$ws.Range("H67").Value = 91
$ws.Range("I67").Value = 16
$ws.Range("J67").Value = 15
$ws.Range("K67").Value = 5
$ws.Range("L67").Value = 7

$ws.Range("B68").Value = 28516
$ws.Range("C68").Value = 28065
$ws.Range("D68").Value = 0.016
$ws.Range("E68").Value = 3775
$ws.Range("F68").Value = 3776
$ws.Range("G68").Value = 22125
$ws.Range("H68").Value = 21645
$ws.Range("I68").Value = 235
$ws.Range("J68").Value = 263
$ws.Range("K68").Value = 2381
$ws.Range("L68").Value = 2381
